# Applies the "Handles float input without breaking stuff" edit:
#   - Marksheet summary rows (10/11/12) get real numbers instead of the
#     placeholder zeros / "Absent" text (student scored 14 right / 1 wrong /
#     13 not-attempted out of 28, total 55 out of 112).
#   - The label cells in column A for rows 10/11/12 pick up the bold
#     "mtitleStyle" formatting (style index 4) that the rest of the header
#     block already uses.
#   - The third question block (columns G/H, "Student Ans"/"Correct Ans")
#     is dropped entirely, and the second block (D/E) collapses down to
#     just the first two data rows.
#   - Column A ("Student Ans") gets filled in with the options the student
#     actually answered; matches are styled with the green "correctStyle"
#     and mismatches with the red "incorrectStyle".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-FormatFrom($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# ---- Row 10 / 11 / 12 labels pick up the "mtitleStyle" (style index 4) ----
Set-FormatFrom "A9" "A10"
Set-FormatFrom "A9" "A11"
Set-FormatFrom "A9" "A12"

# ---- Row 10 (No.) : Right / Wrong / Not Attempt / Max ----
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 28

# ---- Row 11 (Marking) ----
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# ---- Row 12 (Total) ----
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "55/112"

# ---- Row 15 : drop the third "Student Ans"/"Correct Ans" header pair ----
$ws.Range("G15").Clear() | Out-Null
$ws.Range("H15").Clear() | Out-Null

# ---- Row 16 : second block keeps D/E, gets the (correct) student answer ----
Set-FormatFrom "B10" "D16"
$ws.Range("D16").Value = "Option A"
$ws.Range("G16").Clear() | Out-Null
$ws.Range("H16").Clear() | Out-Null

# ---- Row 17 : second block keeps D/E, gets the (incorrect) student answer ----
Set-FormatFrom "C10" "D17"
$ws.Range("D17").Value = "Option B"
$ws.Range("G17").Clear() | Out-Null
$ws.Range("H17").Clear() | Out-Null

# ---- Row 18 : first block gets a correct student answer; third block dropped ----
Set-FormatFrom "B10" "A18"
$ws.Range("A18").Value = "Option B"
$ws.Range("G18").Clear() | Out-Null
$ws.Range("H18").Clear() | Out-Null

# ---- Row 19 : first block correct; second & third blocks dropped ----
Set-FormatFrom "B10" "A19"
$ws.Range("A19").Value = "Option C"
$ws.Range("D19").Clear() | Out-Null
$ws.Range("E19").Clear() | Out-Null
$ws.Range("G19").Clear() | Out-Null
$ws.Range("H19").Clear() | Out-Null

# ---- Row 20 : not attempted; second & third blocks dropped ----
$ws.Range("D20").Clear() | Out-Null
$ws.Range("E20").Clear() | Out-Null
$ws.Range("G20").Clear() | Out-Null
$ws.Range("H20").Clear() | Out-Null

# ---- Row 21 : first block correct; second & third blocks dropped ----
Set-FormatFrom "B10" "A21"
$ws.Range("A21").Value = "Option C"
$ws.Range("D21").Clear() | Out-Null
$ws.Range("E21").Clear() | Out-Null
$ws.Range("G21").Clear() | Out-Null
$ws.Range("H21").Clear() | Out-Null

# ---- Row 22 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A22"
$ws.Range("A22").Value = "Option D"
$ws.Range("D22").Clear() | Out-Null
$ws.Range("E22").Clear() | Out-Null

# ---- Row 23 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A23"
$ws.Range("A23").Value = "Option D"
$ws.Range("D23").Clear() | Out-Null
$ws.Range("E23").Clear() | Out-Null

# ---- Row 24 : not attempted; second block dropped ----
$ws.Range("D24").Clear() | Out-Null
$ws.Range("E24").Clear() | Out-Null

# ---- Row 25 : not attempted; second block dropped ----
$ws.Range("D25").Clear() | Out-Null
$ws.Range("E25").Clear() | Out-Null

# ---- Row 26 : not attempted; second block dropped ----
$ws.Range("D26").Clear() | Out-Null
$ws.Range("E26").Clear() | Out-Null

# ---- Row 27 : not attempted; second block dropped ----
$ws.Range("D27").Clear() | Out-Null
$ws.Range("E27").Clear() | Out-Null

# ---- Row 28 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A28"
$ws.Range("A28").Value = "Option D"
$ws.Range("D28").Clear() | Out-Null
$ws.Range("E28").Clear() | Out-Null

# ---- Row 29 : not attempted; second block dropped ----
$ws.Range("D29").Clear() | Out-Null
$ws.Range("E29").Clear() | Out-Null

# ---- Row 30 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A30"
$ws.Range("A30").Value = "Option B"
$ws.Range("D30").Clear() | Out-Null
$ws.Range("E30").Clear() | Out-Null

# ---- Row 31 : not attempted; second block dropped ----
$ws.Range("D31").Clear() | Out-Null
$ws.Range("E31").Clear() | Out-Null

# ---- Row 32 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A32"
$ws.Range("A32").Value = "Option C"
$ws.Range("D32").Clear() | Out-Null
$ws.Range("E32").Clear() | Out-Null

# ---- Row 33 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A33"
$ws.Range("A33").Value = "Option D"
$ws.Range("D33").Clear() | Out-Null
$ws.Range("E33").Clear() | Out-Null

# ---- Row 34 : not attempted; second block dropped ----
$ws.Range("D34").Clear() | Out-Null
$ws.Range("E34").Clear() | Out-Null

# ---- Row 35 : not attempted; second block dropped ----
$ws.Range("D35").Clear() | Out-Null
$ws.Range("E35").Clear() | Out-Null

# ---- Row 36 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A36"
$ws.Range("A36").Value = "Option A"
$ws.Range("D36").Clear() | Out-Null
$ws.Range("E36").Clear() | Out-Null

# ---- Row 37 : not attempted; second block dropped ----
$ws.Range("D37").Clear() | Out-Null
$ws.Range("E37").Clear() | Out-Null

# ---- Row 38 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A38"
$ws.Range("A38").Value = "Option A"
$ws.Range("D38").Clear() | Out-Null
$ws.Range("E38").Clear() | Out-Null

# ---- Row 39 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A39"
$ws.Range("A39").Value = "Option D"
$ws.Range("D39").Clear() | Out-Null
$ws.Range("E39").Clear() | Out-Null

# ---- Row 40 : first block correct; second block dropped ----
Set-FormatFrom "B10" "A40"
$ws.Range("A40").Value = "Option D"
$ws.Range("D40").Clear() | Out-Null
$ws.Range("E40").Clear() | Out-Null
